$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The per-row tuples (Id, Antal, Ost, Nord) in rows 12-21 get reshuffled
# among those same rows (all other columns in each row are unchanged).
$rows = @(
    @{ Row = 12; A = 111378933; I = "25";  Q = 505597.6535686332; R = 6913018.009825628 },
    @{ Row = 13; A = 111378946; I = "100"; Q = 505602.791734456;  R = 6913005.013642685 },
    @{ Row = 14; A = 111378913; I = "25";  Q = 505607.407264018;  R = 6913026.386327411 },
    @{ Row = 15; A = 111378964; I = "5";   Q = 505627.1571942444; R = 6912898.692122459 },
    @{ Row = 16; A = 111378884; I = "50";  Q = 505596.2310213979; R = 6913034.263345711 },
    @{ Row = 17; A = 111378856; I = "10";  Q = 505494.3524330241; R = 6913043.848162009 },
    @{ Row = 18; A = 111378874; I = "50";  Q = 505592.4968292552; R = 6913042.152801346 },
    @{ Row = 19; A = 111378954; I = "15";  Q = 505590.6913760683; R = 6913009.17353364 },
    @{ Row = 20; A = 111378893; I = "25";  Q = 505612.5119866763; R = 6913033.361683531 },
    @{ Row = 21; A = 111378866; I = "10";  Q = 505492.5216403615; R = 6913025.731493607 }
)

foreach ($r in $rows) {
    $ws.Range("A$($r.Row)").Value = $r.A

    # Column I (Antal) is stored as plain text in this sheet (e.g. "25"),
    # not as a number. Temporarily mark the cell as Text so the numeric-
    # looking string isn't auto-converted to a number, then restore the
    # default "Normal" style so no stray number-format metadata is left
    # behind on the cell.
    $ws.Range("I$($r.Row)").NumberFormat = "@"
    $ws.Range("I$($r.Row)").Value = $r.I
    $ws.Range("I$($r.Row)").Style = "Normal"

    $ws.Range("Q$($r.Row)").Value = $r.Q
    $ws.Range("R$($r.Row)").Value = $r.R
}
